{"js": "// Remove the two \"Uttrekket inneholder ETT/EN ...\" + following \"AND/OR\" blocks\n// that were cut from chapter 3.1.3 (the single-archive / single-archive-part\n// output example, and the trailing \"AND/OR\" + \"Arkivdelstatus er satt til ...\"\n// paragraph right after the results table).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Collect the *positions* (not paragraph objects, since this shim's\n// Paragraph has no stable \"index\" property) that must be deleted, matched\n// by their text content so the script is resilient to exact index drift.\nconst positionsToDelete = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n\n  // Block 1: paragraph right after \"Output\" heading that starts the\n  // \"ETT arkiv / EN arkivdel\" example, the blank paragraph after it, and\n  // the following \"AND/OR\" paragraph.\n  if (text.indexOf(\"Uttrekket inneholder\") === 0 && text.indexOf(\"ETT\") !== -1) {\n    positionsToDelete.push(i);\n    if (i + 1 < items.length && items[i + 1].text.trim() === \"\") {\n      positionsToDelete.push(i + 1);\n    }\n    if (i + 2 < items.length && items[i + 2].text.trim() === \"AND/OR\") {\n      positionsToDelete.push(i + 2);\n    }\n  }\n\n  // Block 2: the \"AND/OR\" paragraph that immediately follows the results\n  // table, the blank paragraph after it, and the \"Arkivdelstatus er satt\n  // til ...\" paragraph that explains a status mismatch.\n  if (\n    text.trim() === \"AND/OR\" &&\n    i + 2 < items.length &&\n    items[i + 2].text.indexOf(\"Arkivdelstatus\") === 0\n  ) {\n    positionsToDelete.push(i);\n    if (items[i + 1].text.trim() === \"\") {\n      positionsToDelete.push(i + 1);\n    }\n    positionsToDelete.push(i + 2);\n  }\n}\n\n// Delete from the end of the document backwards so earlier paragraphs'\n// positions stay valid while we delete.\nconst uniquePositions = Array.from(new Set(positionsToDelete)).sort((a, b) => b - a);\nfor (const pos of uniquePositions) {\n  items[pos].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the two \"Uttrekket inneholder ETT/EN ...\" + following \"AND/OR\" blocks\n# that were cut from chapter 3.1.3 (the single-archive / single-archive-part\n# output example, and the trailing \"AND/OR\" + \"Arkivdelstatus er satt til ...\"\n# paragraph right after the results table).\n#\n# Uses Find (character-offset based) instead of Paragraphs.Item(N) indexing so\n# the script is not sensitive to how the host enumerates paragraph-like marks\n# around the table boundary.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaRangeAt($pos) {\n    $r = $d.Range($pos, $pos)\n    $r.Expand(4)  # wdParagraph\n    return $r\n}\n\n# ---- Block 1: \"Uttrekket inneholder ETT arkiv og EN arkivdeler ... avsluttet.\"\n#      paragraph (right after the \"Output\" heading) plus the \"AND/OR\"\n#      paragraph that immediately follows it. ----\n$f1 = $d.Content\n$f1.Find.ClearFormatting()\n$ok1 = $f1.Find.Execute(\"Uttrekket inneholder\")\nif ($ok1) {\n    $block1Start = Get-ParaRangeAt($f1.Start)\n\n    $f1b = $d.Range($block1Start.End, $d.Content.End)\n    $f1b.Find.ClearFormatting()\n    $ok1b = $f1b.Find.Execute(\"AND/OR\")\n    if ($ok1b) {\n        $block1End = Get-ParaRangeAt($f1b.Start)\n        $d.Range($block1Start.Start, $block1End.End).Delete()\n    }\n}\n\n# ---- Block 2: the \"AND/OR\" paragraph right after the results table, plus the\n#      following \"Arkivdelstatus er satt til \u00abARKIVDELSTATUS\u00bb ...\" paragraph. ----\n$f2 = $d.Content\n$f2.Find.ClearFormatting()\n$ok2 = $f2.Find.Execute(\"Dokumentmedium\")\nif ($ok2) {\n    $f2b = $d.Range($f2.End, $d.Content.End)\n    $f2b.Find.ClearFormatting()\n    $ok2b = $f2b.Find.Execute(\"AND/OR\")\n    if ($ok2b) {\n        $block2Start = Get-ParaRangeAt($f2b.Start)\n\n        $f2c = $d.Range($block2Start.End, $d.Content.End)\n        $f2c.Find.ClearFormatting()\n        $ok2c = $f2c.Find.Execute(\"Arkivdelstatus er satt til\")\n        if ($ok2c) {\n            $block2End = Get-ParaRangeAt($f2c.Start)\n            $d.Range($block2Start.Start, $block2End.End).Delete()\n        }\n    }\n}\n"}
